$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "/suppliers/anh1"
$ws.Range("D3").Value = "/suppliers/anh2"
$ws.Range("D4").Value = "/suppliers/anh3"
$ws.Range("D5").Value = "/suppliers/anh4"

$ws.Range("D2").Select()
